$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay text so numeric-looking strings
# like "1.010" or "0.000008884" are not reinterpreted/reformatted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.054.77"
$ws.Range("E2").Value = "  +2.98%  "
$ws.Range("D3").Value = "1.880.68"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -1.12%  "
$ws.Range("D5").Value = "313.99"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "0.4863"
$ws.Range("E7").Value = "  +1.48%  "
$ws.Range("D8").Value = "0.3818"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").Value = "0.07373"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").Value = "21.03"
$ws.Range("E11").Value = "  +3.70%  "
$ws.Range("D12").Value = "0.07792"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "1.906.15"
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").Value = "5.534"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").Value = "6.617"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "91.82"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "0.000008884"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").Value = "1.011"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "28.057.27"
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").Value = "5.128"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "2.167.09"
$ws.Range("E23").Value = "  +3.63%  "
$ws.Range("D24").Value = "10.96"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").Value = "157.31"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "1.927"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").Value = "18.58"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").Value = "2.062"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "4.974"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").Value = "0.08899"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "3.329"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "1.238"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").Value = "0.7748"
$ws.Range("E34").Value = "  +4.64%  "
$ws.Range("D35").Value = "4.660"
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("D36").Value = "2.743"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").Value = "1.126"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").Value = "0.5604"
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("D40").Value = "0.05378"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").Value = "3.009"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "7.079"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").Value = "8.567"
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("D44").Value = "0.1527"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "0.4909"
$ws.Range("E45").Value = "  +2.38%  "
$ws.Range("D46").Value = "10.75"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").Value = "105.76"
$ws.Range("E47").Value = "  +2.97%  "
$ws.Range("D48").Value = "1.009"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").Value = "1.673"
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").Value = "68.57"
$ws.Range("E50").Value = "  +3.24%  "
$ws.Range("D51").Value = "0.06113"
$ws.Range("E51").Value = "  +0.58%  "
